# Timesheet January 2020.xlsx - add the Feb 24 / Feb 25 2020 entries
# (commit: "Added timesheet of date 25-02-2020")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 191: blank separator row (same highlighted style as the other
# day-separator rows, e.g. row 183) - copy formats only.
# ---------------------------------------------------------------------
$ws.Range("A183:C183").Copy()
$ws.Range("A191:C191").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Feb 24 2020 entries (rows 192-197)
# ---------------------------------------------------------------------
$ws.Range("A192").Value = "Feb 24 10:00 to 11:00"
$ws.Range("B192").Value = "Build django backend. Predictions are fetched."
$ws.Range("C192").Value = "Infimetrics"

$ws.Range("A193").Value = "Feb 24 11:00 to 12:00"
$ws.Range("B193").Value = "Modified some code of ml phases. Making some features in django backend."
$ws.Range("C193").Value = "Infimetrics"

$ws.Range("A194").Value = "Feb 24 12:00 to 13:00"
$ws.Range("B194").Value = "Modified some code of ml phases."
$ws.Range("C194").Value = "Infimetrics"

$ws.Range("A195").Value = "Feb 24 13:00 to 14:00"
$ws.Range("B195").Value = "Lunch"
$ws.Range("C195").Value = "Infimetrics"

$ws.Range("A196").Value = "Feb 24 14:00 to 17:36"
$ws.Range("B196").Value = "Embros technology"
$ws.Range("C196").Value = "Chatrapati chowk"

$ws.Range("A197").Value = "Feb 24 17:36 to 18:00"
$ws.Range("B197").Value = "Modified django backend code, got error, working on code."
$ws.Range("C197").Value = "Infimetrics"

# ---------------------------------------------------------------------
# Row 198: blank separator row, same as row 191/183.
# ---------------------------------------------------------------------
$ws.Range("A183:C183").Copy()
$ws.Range("A198:C198").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Feb 25 2020 entries (rows 199-208)
# ---------------------------------------------------------------------
$ws.Range("A199").Value = "Feb 25 10:00 to 11:00"
$ws.Range("B199").Value = "Used HTTP PATCH method for sending multiple data to api server. Found 60-total alarm`ntime + error relationship of sample predictive model."
$ws.Range("C199").Value = "Infimetrics"
# This task note wraps onto two lines, so (like the other multi-line
# notes in the sheet, e.g. row 182/184/190) give it the wrap-text style
# and a taller row to fit it.
$ws.Range("B199").WrapText = $true
$ws.Range("B199").HorizontalAlignment = -4131
$ws.Range("B199").VerticalAlignment = -4108
$ws.Rows.Item(199).RowHeight = 45

$ws.Range("A200").Value = "Feb 25 11:00 to 12:00"
$ws.Range("B200").Value = "Working on integrating ml code in backend."
$ws.Range("C200").Value = "Infimetrics"

$ws.Range("A201").Value = "Feb 25 12:00 to 13:00"
$ws.Range("B201").Value = "Created uml of data files"
$ws.Range("C201").Value = "Infimetrics"

$ws.Range("A202").Value = "Feb 25 13:00 to 13:30"
$ws.Range("A203").Value = "Feb 25 13:30 to 14:00"
$ws.Range("B202").Value = "Created uml for combined data"
$ws.Range("C202").Value = "Infimetrics"
$ws.Range("B203").Value = "Lunch"
$ws.Range("C203").Value = "Infimetrics"

$ws.Range("A204").Value = "Feb 25 14:00 to 15:00"
$ws.Range("B204").Value = "Tried to upload Djangomlapi to heroku, app was not uploaded but not worked."
$ws.Range("C204").Value = "Infimetrics"

$ws.Range("A205").Value = "Feb 25 15:00 to 16:00"
$ws.Range("B205").Value = "Making a normal api with django api"
$ws.Range("C205").Value = "Infimetrics"

$ws.Range("A206").Value = "Feb 25 16:00 to 17:00"
$ws.Range("B206").Value = "Made a normal api using views only, data sent to post gets lost, working on issue"
$ws.Range("C206").Value = "Infimetrics"

$ws.Range("A207").Value = "Feb 25 17:00 to 18:00"
$ws.Range("B207").Value = "working on multivariate timeseries algo."
$ws.Range("C207").Value = "Infimetrics"

$ws.Range("A208").Value = "Feb 25 18:00 to 19:00"
$ws.Range("B208").Value = "Worked on MVTA, tyring some exaples."
$ws.Range("C208").Value = "Infimetrics"

# ---------------------------------------------------------------------
# Scroll / selection, matching where Excel would land after typing the
# last entry.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 192
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A209").Select()
